$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.753.31'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").Value = '1.732.08'
$ws.Range("E3").Value = '  +0.62%  '

$ws.Range("D4").Value = '0.9976'
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = '242.25'
$ws.Range("E5").Value = '  -0.84%  '

$ws.Range("D6").Value = '0.9980'
$ws.Range("E6").Value = '  -0.29%  '

$ws.Range("D7").Value = '0.4918'
$ws.Range("E7").Value = '  +0.88%  '

$ws.Range("D8").Value = '0.2624'
$ws.Range("E8").Value = '  +0.57%  '

$ws.Range("D9").Value = '0.06221'
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").Value = '1.732.23'
$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("D11").Value = '15.99'
$ws.Range("E11").Value = '  +3.46%  '

$ws.Range("D12").Value = '0.06976'
$ws.Range("E12").Value = '  -0.64%  '

$ws.Range("D13").Value = '0.6124'
$ws.Range("E13").Value = '  +2.62%  '

$ws.Range("D14").Value = '4.512'
$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("D15").Value = '77.25'
$ws.Range("E15").Value = '  +0.06%  '

$ws.Range("D16").Value = '0.9983'

$ws.Range("D17").Value = '26.533.87'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '0.9980'
$ws.Range("E18").Value = '  -0.25%  '

$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D21").Value = '1.952.08'
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").Value = '4.468'
$ws.Range("E22").Value = '  -0.29%  '

$ws.Range("D23").Value = '8.572'

$ws.Range("D24").Value = '5.112'
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").Value = '137.90'
$ws.Range("E25").Value = '  +0.37%  '

$ws.Range("D26").Value = '15.35'
$ws.Range("E26").Value = '  +0.73%  '

$ws.Range("E27").Value = '  +3.28%  '

$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '106.51'
$ws.Range("E28").Value = '  -0.49%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.383'
$ws.Range("E29").Value = '  -2.48%  '

$ws.Range("D30").Value = '3.938'
$ws.Range("E30").Value = '  -0.51%  '

$ws.Range("D31").Value = '0.07983'
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").Value = '3.677'
$ws.Range("E32").Value = '  -0.09%  '

$ws.Range("D33").Value = '0.04487'
$ws.Range("E33").Value = '  -0.83%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.607'
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.004'
$ws.Range("E35").Value = '  +1.07%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6227'
$ws.Range("E36").Value = '  +0.30%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9432'
$ws.Range("E37").Value = '  +4.18%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.052'
$ws.Range("E38").Value = '  +4.02%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.422'
$ws.Range("E39").Value = '  +1.86%  '

$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '0.9979'
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01511'
$ws.Range("E41").Value = '  +1.73%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.577'
$ws.Range("E42").Value = '  +3.55%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '99.54'
$ws.Range("E43").Value = '  -1.01%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3862'
$ws.Range("E44").Value = '  +0.41%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.910'
$ws.Range("E45").Value = '  +2.70%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1161'
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05383'
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '7.865'
$ws.Range("E48").Value = '  +2.23%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.30'
$ws.Range("E49").Value = '  +0.66%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '51.71'
$ws.Range("E50").Value = '  +1.43%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.236'
$ws.Range("E51").Value = '  -0.28%  '
